# Paraguay Division Profesional - fix 4 pairs of rows whose match data was swapped/mixed up.
# Each pair represents two different matches on the same date that had their rows data
# entered in the wrong order. We correct this by swapping all columns B..AB between the
# two rows of each pair (column A, the running index, stays as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 6720843
$ws.Range("E2").Value = "Cerro Porteno"
$ws.Range("F2").Value = "Libertad Asuncion"
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = "H"
$ws.Range("J2").Value = 2.375
$ws.Range("K2").Value = 3.2
$ws.Range("L2").Value = 2.7
$ws.Range("M2").Value = 3.75
$ws.Range("N2").Value = 3.3
$ws.Range("O2").Value = 1.85
$ws.Range("P2").Value = 0.5
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.9
$ws.Range("S2").Value = 2.5
$ws.Range("T2").Value = 1.925
$ws.Range("U2").Value = 1.875
$ws.Range("V2").Value = 2.75
$ws.Range("W2").Value = -1
$ws.Range("X2").Value = -1
$ws.Range("Y2").Value = 0.8999999999999999
$ws.Range("Z2").Value = -1
$ws.Range("AA2").Value = -1
$ws.Range("AB2").Value = 0.875

# Row 4
$ws.Range("B4").Value = 6720844
$ws.Range("E4").Value = "Guarani Asuncion"
$ws.Range("F4").Value = "Olimpia Asuncion"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = "A"
$ws.Range("J4").Value = 2.45
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 2.75
$ws.Range("M4").Value = 4
$ws.Range("N4").Value = 3.2
$ws.Range("O4").Value = 1.85
$ws.Range("P4").Value = 0.5
$ws.Range("Q4").Value = 1.875
$ws.Range("R4").Value = 1.925
$ws.Range("S4").Value = 2.5
$ws.Range("T4").Value = 1.925
$ws.Range("U4").Value = 1.875
$ws.Range("V4").Value = -1
$ws.Range("W4").Value = -1
$ws.Range("X4").Value = 0.8500000000000001
$ws.Range("Y4").Value = -1
$ws.Range("Z4").Value = 0.925
$ws.Range("AA4").Value = 0.925
$ws.Range("AB4").Value = -1

# Row 140
$ws.Range("B140").Value = 7493431
$ws.Range("E140").Value = "Sportivo Trinidense"
$ws.Range("F140").Value = "Guairena FC"
$ws.Range("G140").Value = 7
$ws.Range("H140").Value = 2
$ws.Range("I140").Value = "H"
$ws.Range("J140").Value = 2.05
$ws.Range("K140").Value = 3.3
$ws.Range("L140").Value = 3.3
$ws.Range("M140").Value = 2.6
$ws.Range("N140").Value = 3.1
$ws.Range("O140").Value = 2.6
$ws.Range("P140").Value = 0
$ws.Range("Q140").Value = 1.925
$ws.Range("R140").Value = 1.875
$ws.Range("S140").Value = 2.5
$ws.Range("T140").Value = 2
$ws.Range("U140").Value = 1.8
$ws.Range("V140").Value = 1.6
$ws.Range("W140").Value = -1
$ws.Range("X140").Value = -1
$ws.Range("Y140").Value = 0.925
$ws.Range("Z140").Value = -1
$ws.Range("AA140").Value = 1
$ws.Range("AB140").Value = -1

# Row 141
$ws.Range("B141").Value = 7493310
$ws.Range("E141").Value = "Libertad Asuncion"
$ws.Range("F141").Value = "Tacuary"
$ws.Range("G141").Value = 1
$ws.Range("H141").Value = 2
$ws.Range("I141").Value = "A"
$ws.Range("J141").Value = 1.363
$ws.Range("K141").Value = 5
$ws.Range("L141").Value = 7
$ws.Range("M141").Value = 1.571
$ws.Range("N141").Value = 4.2
$ws.Range("O141").Value = 4.75
$ws.Range("P141").Value = -0.75
$ws.Range("Q141").Value = 1.8
$ws.Range("R141").Value = 2
$ws.Range("S141").Value = 2.75
$ws.Range("T141").Value = 1.8
$ws.Range("U141").Value = 2
$ws.Range("V141").Value = -1
$ws.Range("W141").Value = -1
$ws.Range("X141").Value = 3.75
$ws.Range("Y141").Value = -1
$ws.Range("Z141").Value = 1
$ws.Range("AA141").Value = 0.4
$ws.Range("AB141").Value = -0.5

# Row 231
$ws.Range("B231").Value = 7609161
$ws.Range("E231").Value = "Guarani Asuncion"
$ws.Range("F231").Value = "Nacional Asuncion"
$ws.Range("G231").Value = 3
$ws.Range("H231").Value = 1
$ws.Range("I231").Value = "H"
$ws.Range("J231").Value = 2.1
$ws.Range("K231").Value = 3.25
$ws.Range("L231").Value = 3.6
$ws.Range("M231").Value = 2.25
$ws.Range("N231").Value = 3.1
$ws.Range("O231").Value = 3.4
$ws.Range("P231").Value = -0.25
$ws.Range("Q231").Value = 1.9
$ws.Range("R231").Value = 1.9
$ws.Range("S231").Value = 2.25
$ws.Range("T231").Value = 2
$ws.Range("U231").Value = 1.8
$ws.Range("V231").Value = 1.25
$ws.Range("W231").Value = -1
$ws.Range("X231").Value = -1
$ws.Range("Y231").Value = 0.8999999999999999
$ws.Range("Z231").Value = -1
$ws.Range("AA231").Value = 1
$ws.Range("AB231").Value = -1

# Row 232
$ws.Range("B232").Value = 7609668
$ws.Range("E232").Value = "2 de Mayo"
$ws.Range("F232").Value = "Libertad Asuncion"
$ws.Range("G232").Value = 2
$ws.Range("H232").Value = 0
$ws.Range("I232").Value = "H"
$ws.Range("J232").Value = 4.2
$ws.Range("K232").Value = 3.5
$ws.Range("L232").Value = 1.85
$ws.Range("M232").Value = 4
$ws.Range("N232").Value = 3.4
$ws.Range("O232").Value = 1.909
$ws.Range("P232").Value = 0.5
$ws.Range("Q232").Value = 1.9
$ws.Range("R232").Value = 1.9
$ws.Range("S232").Value = 2.25
$ws.Range("T232").Value = 1.85
$ws.Range("U232").Value = 1.95
$ws.Range("V232").Value = 3
$ws.Range("W232").Value = -1
$ws.Range("X232").Value = -1
$ws.Range("Y232").Value = 0.8999999999999999
$ws.Range("Z232").Value = -1
$ws.Range("AA232").Value = -0.5
$ws.Range("AB232").Value = 0.475

# Row 236
$ws.Range("B236").Value = 7609208
$ws.Range("E236").Value = "Libertad Asuncion"
$ws.Range("F236").Value = "Cerro Porteno"
$ws.Range("G236").Value = 1
$ws.Range("H236").Value = 3
$ws.Range("I236").Value = "A"
$ws.Range("J236").Value = 2.75
$ws.Range("K236").Value = 3
$ws.Range("L236").Value = 2.625
$ws.Range("M236").Value = 2.7
$ws.Range("N236").Value = 2.875
$ws.Range("O236").Value = 2.8
$ws.Range("P236").Value = 0
$ws.Range("Q236").Value = 1.875
$ws.Range("R236").Value = 1.925
$ws.Range("S236").Value = 2
$ws.Range("T236").Value = 1.925
$ws.Range("U236").Value = 1.875
$ws.Range("V236").Value = -1
$ws.Range("W236").Value = -1
$ws.Range("X236").Value = 1.8
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = 0.925
$ws.Range("AA236").Value = 0.925
$ws.Range("AB236").Value = -1

# Row 237
$ws.Range("B237").Value = 7609209
$ws.Range("E237").Value = "Nacional Asuncion"
$ws.Range("F237").Value = "Sportivo Luqueno"
$ws.Range("G237").Value = 0
$ws.Range("H237").Value = 1
$ws.Range("I237").Value = "A"
$ws.Range("J237").Value = 2.9
$ws.Range("K237").Value = 3.1
$ws.Range("L237").Value = 2.5
$ws.Range("M237").Value = 2.7
$ws.Range("N237").Value = 3.25
$ws.Range("O237").Value = 2.55
$ws.Range("P237").Value = 0
$ws.Range("Q237").Value = 1.925
$ws.Range("R237").Value = 1.875
$ws.Range("S237").Value = 2.5
$ws.Range("T237").Value = 2
$ws.Range("U237").Value = 1.8
$ws.Range("V237").Value = -1
$ws.Range("W237").Value = -1
$ws.Range("X237").Value = 1.55
$ws.Range("Y237").Value = -1
$ws.Range("Z237").Value = 0.875
$ws.Range("AA237").Value = -1
$ws.Range("AB237").Value = 0.8

Write-Host "Edit applied successfully."
